$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    4  = 3
    5  = 3
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 2
    11 = 3
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 2
    17 = 3
    18 = 0
    19 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
